$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "model_4_1_0"
$ws.Cells.Item(2, 2).Value = 0.5793153337380459
$ws.Cells.Item(2, 3).Value = -0.2534161772129078
$ws.Cells.Item(2, 4).Value = 0.731265211607491
$ws.Cells.Item(2, 5).Value = 0.649680000174929
$ws.Cells.Item(2, 6).Value = 0.4655739665031433
$ws.Cells.Item(2, 7).Value = 0.1982347369194031
$ws.Cells.Item(2, 8).Value = 0.3978950381278992
$ws.Cells.Item(2, 9).Value = 0.2921925783157349

$ws.Cells.Item(3, 1).Value = "model_4_1_1"
$ws.Cells.Item(3, 2).Value = 0.5934364595711265
$ws.Cells.Item(3, 3).Value = -0.2260355343000739
$ws.Cells.Item(3, 4).Value = 0.58731565059817
$ws.Cells.Item(3, 5).Value = 0.5321768054988798
$ws.Cells.Item(3, 6).Value = 0.4499460160732269
$ws.Cells.Item(3, 7).Value = 0.1939043551683426
$ws.Cells.Item(3, 8).Value = 0.6110301613807678
$ws.Cells.Item(3, 9).Value = 0.3901989161968231

$ws.Cells.Item(4, 1).Value = "model_4_1_2"
$ws.Cells.Item(4, 2).Value = 0.6118501712576001
$ws.Cells.Item(4, 3).Value = 0.1142001699080725
$ws.Cells.Item(4, 4).Value = 0.5044302771241904
$ws.Cells.Item(4, 5).Value = 0.497091773921147
$ws.Cells.Item(4, 6).Value = 0.4295674860477448
$ws.Cells.Item(4, 7).Value = 0.1400941908359528
$ws.Cells.Item(4, 8).Value = 0.7337521910667419
$ws.Cells.Item(4, 9).Value = 0.4194623827934265

$ws.Cells.Item(5, 1).Value = "model_4_1_3"
$ws.Cells.Item(5, 2).Value = 0.6463244758400377
$ws.Cells.Item(5, 3).Value = 0.007055900841264617
$ws.Cells.Item(5, 4).Value = 0.4939481396849659
$ws.Cells.Item(5, 5).Value = 0.4775788538934439
$ws.Cells.Item(5, 6).Value = 0.3914145827293396
$ws.Cells.Item(5, 7).Value = 0.1570396423339844
$ws.Cells.Item(5, 8).Value = 0.7492722868919373
$ws.Cells.Item(5, 9).Value = 0.4357375800609589

$ws.Cells.Item(6, 1).Value = "model_4_1_4"
$ws.Cells.Item(6, 2).Value = 0.6716646938212357
$ws.Cells.Item(6, 3).Value = -0.234963522721598
$ws.Cells.Item(6, 4).Value = 0.4938031997392572
$ws.Cells.Item(6, 5).Value = 0.4531631416300884
$ws.Cells.Item(6, 6).Value = 0.363370418548584
$ws.Cells.Item(6, 7).Value = 0.1953163743019104
$ws.Cells.Item(6, 8).Value = 0.7494868636131287
$ws.Cells.Item(6, 9).Value = 0.4561020731925964

$ws.Cells.Item(7, 1).Value = "model_4_1_5"
$ws.Cells.Item(7, 2).Value = 0.688998071849712
$ws.Cells.Item(7, 3).Value = -1.419340977770501
$ws.Cells.Item(7, 4).Value = 0.4587876969428301
$ws.Cells.Item(7, 5).Value = 0.3050162971266849
$ws.Cells.Item(7, 6).Value = 0.3441874980926514
$ws.Cells.Item(7, 7).Value = 0.3826322555541992
$ws.Cells.Item(7, 8).Value = 0.801331639289856
$ws.Cells.Item(7, 9).Value = 0.5796674489974976

$ws.Cells.Item(8, 1).Value = "model_4_1_8"
$ws.Cells.Item(8, 2).Value = 0.73452099720449
$ws.Cells.Item(8, 3).Value = -2.958703157535435
$ws.Cells.Item(8, 4).Value = 0.4953981076502241
$ws.Cells.Item(8, 5).Value = 0.1810691489896832
$ws.Cells.Item(8, 6).Value = 0.2938070297241211
$ws.Cells.Item(8, 7).Value = 0.6260910034179688
$ws.Cells.Item(8, 8).Value = 0.7471253871917725
$ws.Cells.Item(8, 9).Value = 0.6830485463142395

$ws.Cells.Item(9, 1).Value = "model_4_1_6"
$ws.Cells.Item(9, 2).Value = 0.7444548878956501
$ws.Cells.Item(9, 3).Value = -2.137310780394538
$ws.Cells.Item(9, 4).Value = 0.4924937278674798
$ws.Cells.Item(9, 5).Value = 0.2610987530986471
$ws.Cells.Item(9, 6).Value = 0.2828131318092346
$ws.Cells.Item(9, 7).Value = 0.4961831867694855
$ws.Cells.Item(9, 8).Value = 0.7514257431030273
$ws.Cells.Item(9, 9).Value = 0.6162979006767273

$ws.Cells.Item(10, 1).Value = "model_4_1_7"
$ws.Cells.Item(10, 2).Value = 0.7500051154060945
$ws.Cells.Item(10, 3).Value = -2.321294166027235
$ws.Cells.Item(10, 4).Value = 0.5420464360908998
$ws.Cells.Item(10, 5).Value = 0.2840250922246237
$ws.Cells.Item(10, 6).Value = 0.2766706645488739
$ws.Cells.Item(10, 7).Value = 0.5252811908721924
$ws.Cells.Item(10, 8).Value = 0.6780567765235901
$ws.Cells.Item(10, 9).Value = 0.5971755981445312

$ws.Cells.Item(11, 1).Value = "model_4_1_22"
$ws.Cells.Item(11, 2).Value = 0.7806382882113113
$ws.Cells.Item(11, 3).Value = -1.903541032852587
$ws.Cells.Item(11, 4).Value = 0.4156733573925757
$ws.Cells.Item(11, 5).Value = 0.2203934722322193
$ws.Cells.Item(11, 6).Value = 0.2427687793970108
$ws.Cells.Item(11, 7).Value = 0.4592112302780151
$ws.Cells.Item(11, 8).Value = 0.8651677370071411
$ws.Cells.Item(11, 9).Value = 0.6502490043640137

$ws.Cells.Item(12, 1).Value = "model_4_1_23"
$ws.Cells.Item(12, 2).Value = 0.7824054340221815
$ws.Cells.Item(12, 3).Value = -1.955371415459466
$ws.Cells.Item(12, 4).Value = 0.4283182592961361
$ws.Cells.Item(12, 5).Value = 0.2257524660124592
$ws.Cells.Item(12, 6).Value = 0.2408130615949631
$ws.Cells.Item(12, 7).Value = 0.4674084484577179
$ws.Cells.Item(12, 8).Value = 0.8464453816413879
$ws.Cells.Item(12, 9).Value = 0.6457791924476624

$ws.Cells.Item(13, 1).Value = "model_4_1_24"
$ws.Cells.Item(13, 2).Value = 0.7827036232227498
$ws.Cells.Item(13, 3).Value = -1.981733855791226
$ws.Cells.Item(13, 4).Value = 0.4320094411535197
$ws.Cells.Item(13, 5).Value = 0.2261908852970195
$ws.Cells.Item(13, 6).Value = 0.2404830455780029
$ws.Cells.Item(13, 7).Value = 0.4715777933597565
$ws.Cells.Item(13, 8).Value = 0.8409801721572876
$ws.Cells.Item(13, 9).Value = 0.6454135775566101

$ws.Cells.Item(14, 1).Value = "model_4_1_21"
$ws.Cells.Item(14, 2).Value = 0.783005355554469
$ws.Cells.Item(14, 3).Value = -1.970824664933331
$ws.Cells.Item(14, 4).Value = 0.4364500549024534
$ws.Cells.Item(14, 5).Value = 0.230994529534685
$ws.Cells.Item(14, 6).Value = 0.24014912545681
$ws.Cells.Item(14, 7).Value = 0.4698525071144104
$ws.Cells.Item(14, 8).Value = 0.8344053030014038
$ws.Cells.Item(14, 9).Value = 0.6414069533348083

$ws.Cells.Item(15, 1).Value = "model_4_1_20"
$ws.Cells.Item(15, 2).Value = 0.7838483990099369
$ws.Cells.Item(15, 3).Value = -1.9328242379294
$ws.Cells.Item(15, 4).Value = 0.4384129226106575
$ws.Cells.Item(15, 5).Value = 0.236449320930242
$ws.Cells.Item(15, 6).Value = 0.2392161190509796
$ws.Cells.Item(15, 7).Value = 0.463842511177063
$ws.Cells.Item(15, 8).Value = 0.8314990401268005
$ws.Cells.Item(15, 9).Value = 0.6368573307991028

$ws.Cells.Item(16, 1).Value = "model_4_1_11"
$ws.Cells.Item(16, 2).Value = 0.7841268635623901
$ws.Cells.Item(16, 3).Value = -2.11841709116253
$ws.Cells.Item(16, 4).Value = 0.4857550981270535
$ws.Cells.Item(16, 5).Value = 0.2573657293777287
$ws.Cells.Item(16, 6).Value = 0.2389079481363297
$ws.Cells.Item(16, 7).Value = 0.4931950569152832
$ws.Cells.Item(16, 8).Value = 0.7614030838012695
$ws.Cells.Item(16, 9).Value = 0.6194114685058594

$ws.Cells.Item(17, 1).Value = "model_4_1_12"
$ws.Cells.Item(17, 2).Value = 0.784536422051752
$ws.Cells.Item(17, 3).Value = -2.232379755010871
$ws.Cells.Item(17, 4).Value = 0.4982842984373507
$ws.Cells.Item(17, 5).Value = 0.2563926788329911
$ws.Cells.Item(17, 6).Value = 0.2384546846151352
$ws.Cells.Item(17, 7).Value = 0.5112188458442688
$ws.Cells.Item(17, 8).Value = 0.7428520917892456
$ws.Cells.Item(17, 9).Value = 0.6202231049537659

$ws.Cells.Item(18, 1).Value = "model_4_1_10"
$ws.Cells.Item(18, 2).Value = 0.7848892547126597
$ws.Cells.Item(18, 3).Value = -2.199723449027847
$ws.Cells.Item(18, 4).Value = 0.5045369125803407
$ws.Cells.Item(18, 5).Value = 0.2648952486782169
$ws.Cells.Item(18, 6).Value = 0.2380642145872116
$ws.Cells.Item(18, 7).Value = 0.506054162979126
$ws.Cells.Item(18, 8).Value = 0.7335942983627319
$ws.Cells.Item(18, 9).Value = 0.6131313443183899

$ws.Cells.Item(19, 1).Value = "model_4_1_9"
$ws.Cells.Item(19, 2).Value = 0.7849409168517671
$ws.Cells.Item(19, 3).Value = -2.333645170731309
$ws.Cells.Item(19, 4).Value = 0.5754956982896902
$ws.Cells.Item(19, 5).Value = 0.3107279686110793
$ws.Cells.Item(19, 6).Value = 0.2380070388317108
$ws.Cells.Item(19, 7).Value = 0.5272345542907715
$ws.Cells.Item(19, 8).Value = 0.6285310387611389
$ws.Cells.Item(19, 9).Value = 0.5749034881591797

$ws.Cells.Item(20, 1).Value = "model_4_1_19"
$ws.Cells.Item(20, 2).Value = 0.7853476558508574
$ws.Cells.Item(20, 3).Value = -1.967861070136476
$ws.Cells.Item(20, 4).Value = 0.4511464998161359
$ws.Cells.Item(20, 5).Value = 0.2435693634960236
$ws.Cells.Item(20, 6).Value = 0.2375569045543671
$ws.Cells.Item(20, 7).Value = 0.4693837761878967
$ws.Cells.Item(20, 8).Value = 0.8126453757286072
$ws.Cells.Item(20, 9).Value = 0.6309186816215515

$ws.Cells.Item(21, 1).Value = "model_4_1_18"
$ws.Cells.Item(21, 2).Value = 0.7865479563656217
$ws.Cells.Item(21, 3).Value = -1.988964198252747
$ws.Cells.Item(21, 4).Value = 0.4610499642197532
$ws.Cells.Item(21, 5).Value = 0.2497249330861735
$ws.Cells.Item(21, 6).Value = 0.2362284958362579
$ws.Cells.Item(21, 7).Value = 0.4727213680744171
$ws.Cells.Item(21, 8).Value = 0.7979820966720581
$ws.Cells.Item(21, 9).Value = 0.6257843971252441

$ws.Cells.Item(22, 1).Value = "model_4_1_17"
$ws.Cells.Item(22, 2).Value = 0.7867577947907832
$ws.Cells.Item(22, 3).Value = -1.977391013314463
$ws.Cells.Item(22, 4).Value = 0.4618971149862622
$ws.Cells.Item(22, 5).Value = 0.2515942440052799
$ws.Cells.Item(22, 6).Value = 0.2359962910413742
$ws.Cells.Item(22, 7).Value = 0.470890998840332
$ws.Cells.Item(22, 8).Value = 0.7967277765274048
$ws.Cells.Item(22, 9).Value = 0.6242253184318542

$ws.Cells.Item(23, 1).Value = "model_4_1_14"
$ws.Cells.Item(23, 2).Value = 0.7871362402657578
$ws.Cells.Item(23, 3).Value = -1.892209719231161
$ws.Cells.Item(23, 4).Value = 0.4585359746420322
$ws.Cells.Item(23, 5).Value = 0.2573366356287385
$ws.Cells.Item(23, 6).Value = 0.2355774641036987
$ws.Cells.Item(23, 7).Value = 0.4574190676212311
$ws.Cells.Item(23, 8).Value = 0.8017043471336365
$ws.Cells.Item(23, 9).Value = 0.6194357872009277

$ws.Cells.Item(24, 1).Value = "model_4_1_13"
$ws.Cells.Item(24, 2).Value = 0.7879788052212902
$ws.Cells.Item(24, 3).Value = -2.090189000402007
$ws.Cells.Item(24, 4).Value = 0.4906451133452967
$ws.Cells.Item(24, 5).Value = 0.2642863048164229
$ws.Cells.Item(24, 6).Value = 0.2346449941396713
$ws.Cells.Item(24, 7).Value = 0.4887305796146393
$ws.Cells.Item(24, 8).Value = 0.7541627883911133
$ws.Cells.Item(24, 9).Value = 0.613639235496521

$ws.Cells.Item(25, 1).Value = "model_4_1_15"
$ws.Cells.Item(25, 2).Value = 0.7880221814841407
$ws.Cells.Item(25, 3).Value = -1.944176118800485
$ws.Cells.Item(25, 4).Value = 0.4692465899582517
$ws.Cells.Item(25, 5).Value = 0.2610670873984728
$ws.Cells.Item(25, 6).Value = 0.2345969974994659
$ws.Cells.Item(25, 7).Value = 0.4656378328800201
$ws.Cells.Item(25, 8).Value = 0.7858459949493408
$ws.Cells.Item(25, 9).Value = 0.616324245929718

$ws.Cells.Item(26, 1).Value = "model_4_1_16"
$ws.Cells.Item(26, 2).Value = 0.7886764254752756
$ws.Cells.Item(26, 3).Value = -1.975321432938325
$ws.Cells.Item(26, 4).Value = 0.4732165214130417
$ws.Cells.Item(26, 5).Value = 0.2612580822486555
$ws.Cells.Item(26, 6).Value = 0.2338729202747345
$ws.Cells.Item(26, 7).Value = 0.4705636203289032
$ws.Cells.Item(26, 8).Value = 0.7799680233001709
$ws.Cells.Item(26, 9).Value = 0.6161649823188782
